$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on price cells whose new values would otherwise
# be auto-detected by Excel as numbers (single decimal point), so they stay
# stored as plain text strings exactly like the rest of the Price column.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range("D2").Value = "64.715.65"
$ws.Range("E2").Value = "  -1.52%  "
$ws.Range("D3").Value = "3.115.81"
$ws.Range("E3").Value = "  -7.52%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "566.68"
$ws.Range("E5").Value = "  -2.27%  "
$ws.Range("D6").Value = "167.93"
$ws.Range("E6").Value = "  -5.75%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -3.44%  "
$ws.Range("D9").Value = "3.113.80"
$ws.Range("E9").Value = "  -7.36%  "
$ws.Range("E10").Value = "  -5.94%  "
$ws.Range("E11").Value = "  -5.60%  "
$ws.Range("D12").Value = "0.387"
$ws.Range("E12").Value = "  -5.90%  "
$ws.Range("D13").Value = "3.652.78"
$ws.Range("E13").Value = "  -7.65%  "
$ws.Range("E14").Value = "  +0.99%  "
$ws.Range("D15").Value = "26.65"
$ws.Range("E15").Value = "  -7.59%  "
$ws.Range("D16").Value = "64.601.44"
$ws.Range("E16").Value = "  -1.93%  "
$ws.Range("E17").Value = "  -6.00%  "
$ws.Range("D18").Value = "3.114.17"
$ws.Range("E18").Value = "  -7.57%  "
$ws.Range("D19").Value = "5.64"
$ws.Range("E19").Value = "  -3.27%  "
$ws.Range("D20").Value = "12.62"
$ws.Range("E20").Value = "  -7.14%  "
$ws.Range("D21").Value = "352.43"
$ws.Range("E21").Value = "  -3.30%  "
$ws.Range("D22").Value = "7.14"
$ws.Range("E22").Value = "  -4.46%  "
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("D24").Value = "68.63"
$ws.Range("E24").Value = "  -5.29%  "
$ws.Range("D25").Value = "0.489"
$ws.Range("E25").Value = "  -6.98%  "
$ws.Range("D26").Value = "3.266.58"
$ws.Range("E26").Value = "  -7.70%  "
$ws.Range("E27").Value = "  -7.93%  "
$ws.Range("D28").Value = "9.55"
$ws.Range("E28").Value = "  -1.59%  "
$ws.Range("E29").Value = "  -2.25%  "
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("E32").Value = "  -4.30%  "
$ws.Range("D33").Value = "21.59"
$ws.Range("E33").Value = "  -6.22%  "
$ws.Range("E34").Value = "  -9.34%  "
$ws.Range("D35").Value = "6.52"
$ws.Range("E35").Value = "  -6.05%  "
$ws.Range("D36").Value = "1.17"
$ws.Range("E36").Value = "  -5.13%  "
$ws.Range("D37").Value = "158.16"
$ws.Range("E37").Value = "  -1.60%  "
$ws.Range("D38").Value = "1.41"
$ws.Range("E38").Value = "  -6.58%  "
$ws.Range("E39").Value = "  -3.70%  "
$ws.Range("D40").Value = "26.03"
$ws.Range("E40").Value = "  -3.59%  "
$ws.Range("D41").Value = "1.73"
$ws.Range("E41").Value = "  -1.48%  "
$ws.Range("D42").Value = "2.625.38"
$ws.Range("E42").Value = "  -1.45%  "
$ws.Range("D43").Value = "6.04"
$ws.Range("E43").Value = "  -3.82%  "
$ws.Range("E44").Value = "  -7.65%  "
$ws.Range("D45").Value = "4.12"
$ws.Range("E45").Value = "  -4.47%  "
$ws.Range("D46").Value = "39.28"
$ws.Range("E46").Value = "  -0.75%  "
$ws.Range("E47").Value = "  -3.72%  "
$ws.Range("D48").Value = "23.55"
$ws.Range("E48").Value = "  -3.04%  "
$ws.Range("D49").Value = "316.19"
$ws.Range("E49").Value = "  -5.43%  "
$ws.Range("E50").Value = "  -4.82%  "
$ws.Range("E51").Value = "  -1.89%  "

# Restore the default (no explicit) cell style so only the values changed.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
